# Tie effective values to debug logging and trim advanced UI
#
# 1) On "Menu Mock": the "Show Effective Values" toggle (row 94) and the
#    "Last Trigger" arrow (row 95) are replaced by the former
#    "Quick Test Trigger" (row 97) / "Quick Test Now" (row 98) rows, and the
#    rest of the old CSM Advanced tail (old rows 96-105: Last Trigger Reason,
#    the old Quick Test Trigger/Now duplicates, and the seven "Effective: *"
#    rows) is deleted outright.
# 2) On "Providers": the value-provider rows that only existed to back the
#    removed options (EffectiveBasicProvider .. EffectiveParryProvider,
#    LastTriggerProvider, LastTriggerReasonProvider) are deleted too.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Menu Mock"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Menu Mock")

# Row 94 becomes the old "Quick Test Trigger" row (previously row 97).
$ws.Range("B94").Value = "Quick Test Trigger"
$ws.Range("C94").Value = "Arrow"
$ws.Range("D94").Value = "Basic Kill"
$ws.Range("E94").Value = "Basic Kill | Critical Kill | Dismemberment | Decapitation | Parry | Last Enemy | Last Stand"
$ws.Range("F94").Value = "QuickTestTriggerProvider"
$ws.Range("G94").Value = "Which trigger to simulate"

# Row 95 becomes the old "Quick Test Now" row (previously row 98).
$ws.Range("B95").Value = "Quick Test Now"
$ws.Range("C95").Value = "Toggle"
$ws.Range("D95").Value = "Off"
$ws.Range("E95").Value = ""
$ws.Range("F95").Value = ""
$ws.Range("G95").Value = "Toggle to fire the selected trigger once"

# Delete the now-redundant old rows 96-105 in one shot:
#   96  Last Trigger Reason
#   97  Quick Test Trigger   (duplicate of the content now on row 94)
#   98  Quick Test Now       (duplicate of the content now on row 95)
#   99  Effective: Basic Kill
#   100 Effective: Critical Kill
#   101 Effective: Dismemberment
#   102 Effective: Decapitation
#   103 Effective: Parry
#   104 Effective: Last Enemy
#   105 Effective: Last Stand
$ws.Range("A96:A105").EntireRow.Delete()

# ---------------------------------------------------------------------
# Sheet 2: "Providers"
# ---------------------------------------------------------------------
$ps = $wb.Worksheets.Item("Providers")

# Remove the provider rows that only backed the options removed above:
#   18 EffectiveBasicProvider
#   19 EffectiveCriticalProvider
#   20 EffectiveDecapitationProvider
#   21 EffectiveDismembermentProvider
#   22 EffectiveLastEnemyProvider
#   23 EffectiveLastStandProvider
#   24 EffectiveParryProvider
#   30 LastTriggerProvider
#   31 LastTriggerReasonProvider
# Delete the higher-numbered block first so the lower block's row numbers
# stay valid (deleting low-to-high would shift rows 30-31 out from under us).
$ps.Range("A30:A31").EntireRow.Delete()
$ps.Range("A18:A24").EntireRow.Delete()
